$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Rename sheet "2-data" -> "9-data" and bump its index cell A1: 2 -> 9
# ------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("2-data")
$wsData.Name = "9-data"
$wsData.Range("A1").Value = 9

# ------------------------------------------------------------------
# 2. Move the active tab/selection from "m-map" onto the renamed
#    "9-data" sheet, with B16:B17 highlighted there.
# ------------------------------------------------------------------
$wsData.Select()
$wsData.Range("B16:B17").Select()

# ------------------------------------------------------------------
# 3. Best-effort formatting touch-ups also present in the edit:
#    shrink Sheet4's default column width slightly.
# ------------------------------------------------------------------
$wsFour = $wb.Worksheets.Item("Sheet4")
$wsFour.StandardWidth = 8.375
